$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Existing row 43 (GCNET CP station): the alt_name in column H is updated
#    from the legacy "CP1" tag to the new "CP1_2021" tag (disambiguating it
#    from the freshly logged device added below in row 49).
# ---------------------------------------------------------------------------
$ws.Range("H43").Value = "CP1_2021"

# ---------------------------------------------------------------------------
# 2. Three new device rows are appended at the bottom of the table for the
#    31 May - 31 July 2021 AWS data cleanup:
#       row 49 -> another GCNET CP IMEI (same station as row 43)
#       row 50 -> SWC_U
#       row 51 -> JAR_U
# ---------------------------------------------------------------------------

# Row 49 follows the same "odd" banding/format as the other data rows
# (e.g. row 43), so clone that row's formatting first, then fill in values.
$ws.Range("A43:J43").Copy()
$ws.Range("A49:J49").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A49").Value = 300534062024750
$ws.Range("B49").Value = "GCNET CP"
$ws.Range("B49:C49").Merge()
$ws.Range("D49").Value = "B"
$ws.Range("E49").Value = "NO"
$ws.Range("F49").Value = "ACTIVE"
$ws.Range("G49").Value = "2021-06-02 13:21:24"
$ws.Range("H49").Value = "CP1"
$ws.Range("I49").Value = 1
$ws.Range("J49").Value = "g"
$ws.Rows.Item(49).RowHeight = 16

# Row 51 - JAR_U. Plain, unstyled row. (Filled in before row 50 below so the
# shared-string table ends up in the same append order as the source file.)
$ws.Range("A51").Value = 300534063814490
$ws.Range("B51").Value = "JAR_U"
$ws.Range("H51").Formula = "=B51"
$ws.Range("I51").Value = 1

$ws.Range("J1").Copy()
$ws.Range("J51").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("J51").Value = "g"

# Row 50 - SWC_U. Column B picks up the "network" style (same as column J)
# rather than the usual asset-name formatting.
$ws.Range("J1").Copy()
$ws.Range("B50").PasteSpecial(-4122)       # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A50").Value = 300534063816770
$ws.Range("B50").Value = "SWC_U"
$ws.Range("H50").Formula = "=B50"
$ws.Range("I50").Value = 1

$ws.Range("J1").Copy()
$ws.Range("J50").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("J50").Value = "g"

# ---------------------------------------------------------------------------
# 3. Update the view/selection so the window shows the newly added rows.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J52").Select()
